# Add a new "2020" data column (Q) to the sheet, mirroring the formatting
# of the existing 2019 column (P), and update the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting of column P (rows 3-33, the data block) into column Q
#    cell-by-cell so each row keeps its own look (borders/number format/font).
for ($r = 3; $r -le 33; $r++) {
    $src = $ws.Cells.Item($r, 16)   # column P
    $dst = $ws.Cells.Item($r, 17)   # column Q
    $src.Copy($dst)
}

# 2) Fill in the new 2020 values.
$ws.Range("Q3").Value2  = 2020
$ws.Range("Q4").Value2  = 53.463696812512026
$ws.Range("Q5").Value2  = 46.05204738706685
$ws.Range("Q6").Value2  = 60.76705279190513
$ws.Range("Q7").Value2  = 46.609654277955656
$ws.Range("Q8").Value2  = 39.785591828762811
$ws.Range("Q9").Value2  = 53.693996785869842
$ws.Range("Q10").Value2 = 49.132459991853935
$ws.Range("Q11").Value2 = 42.132308166831223
$ws.Range("Q12").Value2 = 56.225753650646354
$ws.Range("Q13").Value2 = 28.457427087863305
$ws.Range("Q14").Value2 = 20.524708126577082
$ws.Range("Q15").Value2 = 36.325895173845353
$ws.Range("Q16").Value2 = 37.816151622141014
$ws.Range("Q17").Value2 = 29.032279844170926
$ws.Range("Q18").Value2 = 46.928626462141906
$ws.Range("Q19").Value2 = 51.38232216208695
$ws.Range("Q20").Value2 = 45.862881450184311
$ws.Range("Q21").Value2 = 57.0280888993139
$ws.Range("Q22").Value2 = 44.951834666409091
$ws.Range("Q23").Value2 = 38.216466887636237
$ws.Range("Q24").Value2 = 51.83682668469686
$ws.Range("Q25").Value2 = 82.176148450436926
$ws.Range("Q26").Value2 = 66.965035434789911
$ws.Range("Q27").Value2 = 96.931980629894966
$ws.Range("Q28").Value2 = 56.391242440049062
$ws.Range("Q29").Value2 = 50.844030930786069
$ws.Range("Q30").Value2 = 61.300998533028128
$ws.Range("Q31").Value2 = 54.829571415516767
$ws.Range("Q32").Value2 = 58.407045187583961
$ws.Range("Q33").Value2 = 51.452932817170577

# 3) Update the sheet's stored selection/scroll position to match what Excel
#    recorded after this edit (active cell S34, no frozen/scrolled topLeftCell).
$ws.Range("S34").Select()
